# A new scrape run was recorded. Its timestamped price column is inserted
# into the "Suivi" sheet right before the trailing "nom" / "url_produit"
# columns, which shift one column to the right (CS->CT, CT->CU).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at CS; everything from CS onward (nom, url_produit)
# shifts right by one column.
$ws.Columns("CS").Insert()

# Header of the newly inserted column: timestamp of this scrape run.
$ws.Range("CS1").Value = "2026-02-01 04:19:52"

# This run's prices for the product rows that have a price: same value as
# the previous run's price column (CR), copied straight across.
for ($r = 2; $r -le 80; $r++) {
    $ws.Cells.Item($r, 97).Value = $ws.Cells.Item($r, 96).Value()
}
